$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).ClearFormats()
}

Set-TextValue 'D2' '306.78'
Set-TextValue 'E2' '-0.33%'
Set-TextValue 'G2' '3'
Set-TextValue 'D3' '41.26'
Set-TextValue 'E3' '2.34%'
Set-TextValue 'G3' '3'
Set-TextValue 'D4' '5.108'
Set-TextValue 'E4' '2.09%'
Set-TextValue 'G4' '3'
Set-TextValue 'D5' '0.07605'
Set-TextValue 'E5' '-1.27%'
Set-TextValue 'G5' '3'
Set-TextValue 'B6' 'GateToken'
Set-TextValue 'C6' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D6' '4.255'
Set-TextValue 'E6' '-0.07%'
Set-TextValue 'G6' '3'
Set-TextValue 'B7' 'FTXToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D7' '1.619'
Set-TextValue 'E7' '0.17%'
Set-TextValue 'G7' '3'
Set-TextValue 'B8' 'BTSEToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D8' '2.488'
Set-TextValue 'E8' '-2.37%'
Set-TextValue 'G8' '3'
Set-TextValue 'B9' 'MXToken'
Set-TextValue 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.9020'
Set-TextValue 'E9' '0.30%'
Set-TextValue 'G9' '3'
Set-TextValue 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.1114'
Set-TextValue 'E10' '12.29%'
Set-TextValue 'G10' '3'
Set-TextValue 'B11' 'WazirX'
Set-TextValue 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1758'
Set-TextValue 'E11' '1.51%'
Set-TextValue 'G11' '3'
Set-TextValue 'B12' 'MandalaExchangeToken'
Set-TextValue 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.09170'
Set-TextValue 'E12' '3.00%'
Set-TextValue 'G12' '3'
Set-TextValue 'B13' 'BitrueCoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.04291'
Set-TextValue 'E13' '-2.09%'
Set-TextValue 'G13' '3'
Set-TextValue 'B14' 'BitMartToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.1049'
Set-TextValue 'E14' '-0.66%'
Set-TextValue 'G14' '3'
Set-TextValue 'B15' 'BitForexToken'
Set-TextValue 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001255'
Set-TextValue 'E15' '-0.81%'
Set-TextValue 'G15' '3'
Set-TextValue 'B16' 'TigerCash'
Set-TextValue 'C16' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D16' '0.005821'
Set-TextValue 'E16' '-1.51%'
Set-TextValue 'G16' '3'
Set-TextValue 'B17' 'LEO'
Set-TextValue 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D17' '3.361'
Set-TextValue 'E17' '-0.24%'
Set-TextValue 'G17' '3'
Set-TextValue 'E18' '-2.67%'
Set-TextValue 'G18' '3'
Set-TextValue 'D19' '6.577'
Set-TextValue 'E19' '-5.68%'
Set-TextValue 'G19' '3'
Set-TextValue 'D20' '0.1356'
Set-TextValue 'E20' '0.53%'
Set-TextValue 'G20' '3'
Set-TextValue 'D21' '0.2683'
Set-TextValue 'G21' '3'
Set-TextValue 'D22' '0.04186'
Set-TextValue 'E22' '-1.41%'
Set-TextValue 'G22' '3'
Set-TextValue 'E23' '2.05%'
Set-TextValue 'G23' '3'
Set-TextValue 'D24' '0.004075'
Set-TextValue 'E24' '0.05%'
Set-TextValue 'G24' '3'
Set-TextValue 'E25' '6.67%'
Set-TextValue 'G25' '3'
Set-TextValue 'D26' '0.0003008'
Set-TextValue 'E26' '0.91%'
Set-TextValue 'G26' '3'
Set-TextValue 'G27' '3'
Set-TextValue 'G28' '3'
Set-TextValue 'G29' '3'
Set-TextValue 'G30' '3'
Set-TextValue 'G31' '3'
Set-TextValue 'G32' '3'
Set-TextValue 'G33' '3'
Set-TextValue 'G34' '3'
Set-TextValue 'G35' '3'
Set-TextValue 'G36' '3'
Set-TextValue 'G37' '3'
Set-TextValue 'E38' '1.11%'
Set-TextValue 'G38' '3'
Set-TextValue 'D39' '0.05173'
Set-TextValue 'E39' '-0.28%'
Set-TextValue 'G39' '3'
Set-TextValue 'D40' '0.007755'
Set-TextValue 'E40' '-2.51%'
Set-TextValue 'G40' '3'
Set-TextValue 'D41' '0.1296'
Set-TextValue 'E41' '-2.19%'
Set-TextValue 'G41' '3'
Set-TextValue 'D42' '0.006956'
Set-TextValue 'E42' '1.83%'
Set-TextValue 'G42' '3'
Set-TextValue 'D43' '0.001971'
Set-TextValue 'E43' '0.71%'
Set-TextValue 'G43' '3'
Set-TextValue 'D44' '0.008542'
Set-TextValue 'E44' '15.00%'
Set-TextValue 'G44' '3'
Set-TextValue 'E45' '-8.35%'
Set-TextValue 'G45' '3'
Set-TextValue 'D46' '0.00006542'
Set-TextValue 'E46' '-1.45%'
Set-TextValue 'G46' '3'
Set-TextValue 'D47' '0.00000000750'
Set-TextValue 'G47' '3'
Set-TextValue 'B48' 'BOLO'
Set-TextValue 'C48' 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue 'D48' '0.008999'
Set-TextValue 'E48' '187.14%'
Set-TextValue 'G48' '3'
Set-TextValue 'B49' 'CoinbaseStockToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue 'D49' '0.004402'
Set-TextValue 'E49' '-11.94%'
Set-TextValue 'G49' '3'
Set-TextValue 'D50' '0.00002101'
Set-TextValue 'G50' '3'
Set-TextValue 'D51' '0.0002001'
Set-TextValue 'G51' '3'
